# Updates the crypto price/volume table to the latest scrape.
# Column D ("Price") holds plain-text numbers (often with thousands
# separators that collide with the decimal point, e.g. "23.307.82"),
# so cells whose new value *would* auto-parse as a number are force-
# formatted as Text first -- this mirrors formatting a cell as Text in
# the Excel UI before typing a numeric-looking value into it, keeping it
# a literal string instead of silently becoming 302.54000000000002, etc.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.307.82'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '1.623.14'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '302.54'
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3754'
$ws.Range('E7').Value = '  +0.82%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3617'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '51.29'
$ws.Range('E9').Value = '  -1.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08134'
$ws.Range('E10').Value = '  +0.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.219'
$ws.Range('E11').Value = '  -2.34%  '
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.457'
$ws.Range('E14').Value = '  -2.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001234'
$ws.Range('E15').Value = '  -2.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.268'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = '1.618.81'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.95'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06931'
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.47'
$ws.Range('E20').Value = '  -3.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.525'
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('D24').Value = '23.301.24'
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('E25').Value = '  +2.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.072'
$ws.Range('E26').Value = '  +1.70%  '
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '150.51'
$ws.Range('E28').Value = '  -0.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.271'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.65'
$ws.Range('E30').Value = '  -2.25%  '
$ws.Range('D31').Value = '1.797.45'
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.710'
$ws.Range('E32').Value = '  -0.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.119'
$ws.Range('E33').Value = '  -7.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.059'
$ws.Range('E34').Value = '  +11.11%  '
$ws.Range('E35').Value = '  +9.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02741'
$ws.Range('E36').Value = '  -3.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08763'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2474'
$ws.Range('E38').Value = '  -1.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.07086'
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.957'
$ws.Range('E40').Value = '  -1.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6957'
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.325'
$ws.Range('E42').Value = '  -3.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '15.93'
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.99'
$ws.Range('E44').Value = '  -3.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6439'
$ws.Range('E45').Value = '  -0.85%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.260'
$ws.Range('E47').Value = '  -2.70%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.949'
$ws.Range('E48').Value = '  -1.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07956'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.73'
$ws.Range('E50').Value = '  -2.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.182'
$ws.Range('E51').Value = '  -1.13%  '
